$d = $word.ActiveDocument

# Locate the "© 2020 ... Powered by Jekyll ..." footer paragraph. The two
# paragraphs immediately preceding it are a "Ver no Jupiter ..." line and a
# blank spacer paragraph; all three are removed as a block, leaving the
# trailing blank paragraph + page-break paragraph intact.
$copyrightIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Powered by Jekyll*") {
        $copyrightIndex = $i
        break
    }
}

if ($copyrightIndex -ge 1) {
    # Delete from the end backwards so earlier indices stay valid.
    $d.Paragraphs.Item($copyrightIndex).Range.Delete()
    $d.Paragraphs.Item($copyrightIndex - 1).Range.Delete()
    $d.Paragraphs.Item($copyrightIndex - 2).Range.Delete()
}

Write-Output "copyrightIndex=$copyrightIndex"
